# Append a new "Last Updated" log entry row to the Log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRows = $ws.UsedRange.Rows.Count
$newRow = $usedRows + 1

$ws.Cells.Item($newRow, 1).Value = "2025-06-20 14:57:42"

$wb.Save()
